$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: data, hora, preco, site, cor
$ws.Cells.Item(20, 1).Value = 45211
$ws.Cells.Item(20, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20, 2).Value = "16:08"
$ws.Cells.Item(20, 3).Value = 1829
$ws.Cells.Item(20, 4).Value = "amazon"
$ws.Cells.Item(20, 5).Value = "preto"

# Row 21: data, hora, preco, site, cor
$ws.Cells.Item(21, 1).Value = 45211
$ws.Cells.Item(21, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21, 2).Value = "16:09"
$ws.Cells.Item(21, 3).Value = 1829
$ws.Cells.Item(21, 4).Value = "amazon"
$ws.Cells.Item(21, 5).Value = "preto"
